# ResultSearchQuestionCart.xlsx update:
# Cells E5:E8 flip from "PASS" (green "Good" style) to "FAIL" (red "Bad" style).
#
# The sheet already carries a "Good" style (bold dark-green font on a light
# green fill) used for every Pass/Fail cell. We derive the new "Bad" style
# (bold dark-red font on a light red fill) from that existing style so the
# border/alignment/font-size stay identical and only the colors change -
# matching the new font+fill+cellXf entries appended in the workbook's
# styles part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new "Bad" look on a single cell first (E5) so only one extra
# font/fill/cellXf combination gets minted, then propagate that exact
# format to E6:E8 via a format-only copy before writing the new text.
# (RGB() isn't available here, so the BGR-packed long values are given
# directly: 9C0006 -> 393372, FFC7CE -> 13551615)
$e5 = $ws.Range("E5")
$e5.Font.Color = 393372
$e5.Interior.Color = 13551615

$e5.Copy()
$ws.Range("E6:E8").PasteSpecial(-4122)

$ws.Range("E5:E8").Value = "FAIL"
